$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Row 31's "Value" (B31) was entered as 3317, missing the decimal point that
# every other row in the series has (e.g. 3.452, 3.421, ...). Correct the
# typo to 3.317 - this also fixes the downstream "% Change vs Last Year"
# shared formula results in C19 (=(B19/B31-1)*100) and C31
# (=(B31/B43-1)*100), which were wildly wrong because of the stray 1000x.
$ws.Range("B31").Value = 3.317

# Match the "#,##0.000" 3-decimal display format now needed for this cell.
$ws.Range("B31").NumberFormat = "#,##0.000"
